# "Save photos to folder and add new lifts to excel"
# Append 5 new lift records (rows 4-8) to the "lifts" sheet.
#
# Row 4 keeps the same explicit cell style as the existing data rows (2-3),
# while rows 5-8 are written with the plain/default style - matching how
# the source workbook's new rows were authored.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lifts")

# Row 4 - ID 1, styled like the pre-existing rows
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "NONE"
$ws.Range("C4").Value = "S"
$ws.Range("D4").Value = "O"
$ws.Range("E4").Value = "Note"

# Row 5 - ID 2, default style
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "NONE"
$ws.Range("C5").Value = "S"
$ws.Range("D5").Value = "O"
$ws.Range("E5").Value = "Note"
$ws.Range("A5:E5").Style = "Normal"

# Row 6 - ID 3, default style
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "NONE"
$ws.Range("C6").Value = "S"
$ws.Range("D6").Value = "O"
$ws.Range("E6").Value = "Note"
$ws.Range("A6:E6").Style = "Normal"

# Row 7 - ID 4, default style
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "NONE"
$ws.Range("C7").Value = "S"
$ws.Range("D7").Value = "O"
$ws.Range("E7").Value = "Note"
$ws.Range("A7:E7").Style = "Normal"

# Row 8 - ID 5, new site/opening codes, default style
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "NONE"
$ws.Range("C8").Value = "047G"
$ws.Range("D8").Value = "047G-P"
$ws.Range("E8").Value = "Note"
$ws.Range("A8:E8").Style = "Normal"

# Match the author's final selection/active cell
$ws.Range("E4").Select()
